# Community CRUD dummy-data sheet: swap the "ino22" placeholder user id for
# "iin22" and collapse the duplicate ",1,1);" shared string so the SQL
# INSERT statements built by column F use the new id for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newC = ",'iin22','더미용','냉무',null,0,TO_DATE('2021-08-29 00:00:30','YYYY-MM-DD HH24:MI:SS'),'#연봉',0,null,"
$newE = ",1,1);"

for ($r = 1; $r -le 49; $r++) {
    $ws.Cells.Item($r, 3).Value = $newC
    $ws.Cells.Item($r, 5).Value = $newE
}

# Restore the view state recorded for the sheet after the edit.
$ws.Application.ActiveWindow.ScrollRow = 37
$ws.Range("C5").Select()
